$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new value for column O ("Duration")
$durationUpdates = @{
    "9" = 0.001238584518432617
    "10" = 0.0009899139404296875
    "11" = 0.02717375755310059
    "12" = 0.002861261367797852
    "13" = 0.01743292808532715
    "14" = 0.123436450958252
    "15" = 0.003000974655151367
    "16" = 0.006684064865112305
    "17" = 0.04587340354919434
    "18" = 0.002004146575927734
    "19" = 0.02090883255004883
    "20" = 0.0009725093841552734
    "21" = 0.01871132850646973
    "22" = 0.001993656158447266
    "23" = 0
    "24" = 0.0009891986846923828
    "25" = 0.007002830505371094
    "26" = 0.06189751625061035
    "27" = 0.002984762191772461
    "30" = 0.0030059814453125
    "31" = 0.0363919734954834
    "33" = 0.0009946823120117188
    "34" = 0.05583310127258301
    "35" = 0.02160072326660156
    "36" = 0.0008256435394287109
    "37" = 0.08654332160949707
    "38" = 0.04495835304260254
    "39" = 0.06683254241943359
    "40" = 0.04353547096252441
    "41" = 0.01604652404785156
    "42" = 0.02392148971557617
    "43" = 0.0176842212677002
    "44" = 0.03731489181518555
    "45" = 0.002000331878662109
    "46" = 0.03380107879638672
    "47" = 0.002107381820678711
    "48" = 0.01739501953125
    "49" = 0.01351213455200195
    "50" = 0.00199127197265625
    "51" = 0.006040573120117188
    "52" = 0.002001285552978516
    "53" = 0.01565098762512207
    "54" = 0.002007961273193359
    "55" = 0.04849100112915039
    "57" = 0.02603363990783691
    "58" = 0.01462340354919434
    "60" = 0.0492253303527832
    "62" = 0.006050348281860352
    "63" = 0.07803916931152344
    "64" = 0.06633329391479492
    "65" = 0.4735481739044189
    "66" = 0.01341080665588379
    "67" = 0.003023862838745117
    "68" = 0.001993656158447266
    "69" = 0.05604648590087891
    "70" = 0.04452061653137207
    "71" = 16.38140368461609
    "72" = 186.1293351650238
    "73" = 0
    "74" = 0.001004219055175781
    "75" = 0.002202749252319336
    "77" = 0
    "79" = 0.002198457717895508
    "80" = 0.00198674201965332
    "81" = 0.0674440860748291
    "82" = 0
    "83" = 0.002001762390136719
    "84" = 0.5386946201324463
    "85" = 0.00502324104309082
    "86" = 0.00593876838684082
    "87" = 0.1505951881408691
    "88" = 0.1274991035461426
    "89" = 0.004006147384643555
    "90" = 0.419173002243042
    "91" = 0.06857752799987793
    "92" = 0.06133031845092773
    "93" = 0.01563596725463867
    "94" = 0.003249406814575195
    "95" = 0.0103905200958252
    "96" = 2.70869255065918
    "97" = 0.0217583179473877
    "98" = 0.06321144104003906
    "99" = 0.05058979988098145
    "100" = 0.0568697452545166
    "101" = 0.03600049018859863
    "103" = 0.03467202186584473
    "104" = 2.23248291015625
    "108" = 0.04092812538146973
    "109" = 0.003993511199951172
    "110" = 0.05959296226501465
    "111" = 0.05014777183532715
    "112" = 0.04639625549316406
    "113" = 0.0009820461273193359
    "114" = 0.04128050804138184
    "115" = 0.02747178077697754
    "116" = 0.0009164810180664062
    "117" = 0.002113819122314453
    "118" = 0.0009980201721191406
    "119" = 0
    "120" = 0.001997232437133789
    "121" = 0.02913784980773926
    "122" = 0.03264355659484863
    "124" = 0.00244140625
    "125" = 0.006496191024780273
    "126" = 0
    "127" = 0.001000642776489258
    "129" = 0.001044034957885742
    "130" = 0.002027511596679688
    "131" = 0.0009920597076416016
    "132" = 0.009003639221191406
    "133" = 0.001764297485351562
    "134" = 0.00101017951965332
    "135" = 0.002320766448974609
    "137" = 0
    "139" = 0.01157760620117188
    "140" = 0.03242611885070801
    "141" = 0.002999067306518555
    "142" = 0.05230903625488281
    "143" = 0.001664876937866211
    "144" = 0.0202934741973877
    "145" = 0.003092765808105469
    "146" = 0.002726316452026367
    "147" = 0.1100273132324219
    "148" = 0.03312110900878906
    "149" = 0.2554035186767578
    "150" = 0.2010571956634521
    "151" = 4.446112155914307
    "152" = 0.106736421585083
    "153" = 0
    "154" = 0.004574060440063477
    "155" = 0.1401519775390625
    "156" = 0
    "157" = 0.01747679710388184
    "158" = 0.1935446262359619
    "159" = 0.1466183662414551
    "160" = 0.004705429077148438
    "162" = 0.01831221580505371
    "163" = 0.007512092590332031
    "164" = 0.1833932399749756
    "165" = 0.120621919631958
    "166" = 0.0722966194152832
    "167" = 0.1333177089691162
    "169" = 0.01279473304748535
    "170" = 0
    "171" = 0.02595090866088867
    "172" = 0.02843761444091797
    "173" = 0.04254531860351562
    "174" = 0.05129647254943848
    "175" = 0.001169204711914062
    "176" = 0.01341414451599121
    "177" = 0.003011703491210938
}

foreach ($row in $durationUpdates.Keys) {
    $ws.Cells.Item([int]$row, 15).Value = $durationUpdates[$row]
}
